$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '70.115.68'
Set-TextValue 'E2' '  +3.06%  '
Set-TextValue 'D3' '2.454.79'
Set-TextValue 'E3' '  +1.35%  '
Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  -0.02%  '
Set-TextValue 'D5' '567.89'
Set-TextValue 'E5' '  +2.18%  '
Set-TextValue 'D6' '167.66'
Set-TextValue 'E6' '  +4.40%  '
Set-TextValue 'E7' '  -0.05%  '
Set-TextValue 'D8' '0.513'
Set-TextValue 'E8' '  +0.23%  '
Set-TextValue 'D9' '0.175'
Set-TextValue 'E9' '  +12.60%  '
Set-TextValue 'D10' '2.450.97'
Set-TextValue 'E10' '  +1.30%  '
Set-TextValue 'E11' '  -1.48%  '
Set-TextValue 'D12' '0.336'
Set-TextValue 'E12' '  +3.30%  '
Set-TextValue 'E13' '  -1.15%  '
Set-TextValue 'E14' '  +8.91%  '
Set-TextValue 'D15' '70.048.09'
Set-TextValue 'E15' '  +3.09%  '
Set-TextValue 'D16' '2.911.71'
Set-TextValue 'E16' '  +0.20%  '
Set-TextValue 'D17' '24.20'
Set-TextValue 'E17' '  +5.50%  '
Set-TextValue 'D18' '2.460.77'
Set-TextValue 'E18' '  +0.40%  '
Set-TextValue 'B19' 'Chainlink'
Set-TextValue 'C19' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D19' '10.89'
Set-TextValue 'E19' '  +5.89%  '
Set-TextValue 'B20' 'Uniswap'
Set-TextValue 'C20' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D20' '7.24'
Set-TextValue 'E20' '  +6.86%  '
Set-TextValue 'B21' 'BitcoinCash'
Set-TextValue 'C21' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D21' '343.69'
Set-TextValue 'E21' '  +2.85%  '
Set-TextValue 'B22' 'Polkadot'
Set-TextValue 'C22' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D22' '3.90'
Set-TextValue 'E22' '  +3.87%  '
Set-TextValue 'B23' 'Binance-PegBSC-USD'
Set-TextValue 'C23' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 'D23' '1.33'
Set-TextValue 'E23' '  +33.77%  '
Set-TextValue 'E24' '  +7.99%  '
Set-TextValue 'E25' '  -0.02%  '
Set-TextValue 'D26' '66.56'
Set-TextValue 'E26' '  -0.30%  '
Set-TextValue 'E27' '  +7.64%  '
Set-TextValue 'B28' 'Aptos'
Set-TextValue 'C28' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D28' '8.54'
Set-TextValue 'E28' '  +6.84%  '
Set-TextValue 'B29' 'WrappedeETH'
Set-TextValue 'C29' 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue 'D29' '2.580.03'
Set-TextValue 'E29' '  +1.16%  '
Set-TextValue 'D30' '0.0₃0871'
Set-TextValue 'E30' '  +8.87%  '
Set-TextValue 'E31' '  +4.56%  '
Set-TextValue 'E32' '  +11.43%  '
Set-TextValue 'D33' '450.45'
Set-TextValue 'E33' '  +7.08%  '
Set-TextValue 'D34' '1.00'
Set-TextValue 'E34' '  +0.03%  '
Set-TextValue 'E35' '  +2.58%  '
Set-TextValue 'D36' '162.18'
Set-TextValue 'E36' '  +0.85%  '
Set-TextValue 'D37' '19.10'
Set-TextValue 'E37' '  +0.69%  '
Set-TextValue 'E38' '  +7.47%  '
Set-TextValue 'E39' '  -0.03%  '
Set-TextValue 'E40' '  +3.10%  '
Set-TextValue 'E41' '  +4.91%  '
Set-TextValue 'E42' '  +7.01%  '
Set-TextValue 'D43' '4.44'
Set-TextValue 'E43' '  +4.38%  '
Set-TextValue 'E44' '  +5.40%  '
Set-TextValue 'D45' '2.15'
Set-TextValue 'E45' '  +8.98%  '
Set-TextValue 'E46' '  +2.80%  '
Set-TextValue 'D47' '133.53'
Set-TextValue 'E47' '  +4.23%  '
Set-TextValue 'D48' '0.0725'
Set-TextValue 'E48' '  +1.73%  '
Set-TextValue 'D49' '0.492'
Set-TextValue 'E49' '  +4.15%  '
Set-TextValue 'D50' '0.564'
Set-TextValue 'E50' '  +1.86%  '
Set-TextValue 'E51' '  +2.10%  '
